# Applies the data update described by the diff:
# - Inserts two new rows at the top of the data (tv, carregador)
# - Keeps/duplicates the original 4 data rows, with an extra "caneta" row
# - Ends with the same 4-row block repeated again (total 12 data rows + header)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("tv",         7,  "eletrônico",             1000),
    @("carregador", 2,  "eletrônico",             7),
    @("varistor",   2,  "componente eletrônico",  0.35),
    @("lapis",      10, "item de escola",         2),
    @("cola",       50, "item de escola",         2),
    @("lapis",      10, "item de escola",         2.5),
    @("caneta",     80, "item de escola",         3.9),
    @("varistor",   2,  "componente eletrônico",  0.35),
    @("lapis",      10, "item de escola",         2),
    @("cola",       50, "item de escola",         2),
    @("lapis",      10, "item de escola",         2.5),
    @("caneta",     80, "item de escola",         3.9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}
